$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 104, shifting existing rows 104:229 down to 105:230
$ws.Rows.Item(104).Insert()

# Populate the newly inserted row 104 with the new weekly price record
$ws.Cells.Item(104, 1).Value = 3
$ws.Cells.Item(104, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(104, 3).Value = "Coquimbo"
$ws.Cells.Item(104, 4).Value = 44966
$ws.Cells.Item(104, 5).Value = 5
$ws.Cells.Item(104, 6).Value = 100112052
$ws.Cells.Item(104, 7).Value = "Albahaca"
$ws.Cells.Item(104, 8).Value = "Sin especificar"
$ws.Cells.Item(104, 9).Value = "Primera"
$ws.Cells.Item(104, 10).Value = 110
$ws.Cells.Item(104, 11).Value = 5500
$ws.Cells.Item(104, 12).Value = 6000
$ws.Cells.Item(104, 13).Value = 5773
$ws.Cells.Item(104, 14).Value = "$/docena de matas"
$ws.Cells.Item(104, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(104, 16).Value = 962
$ws.Cells.Item(104, 17).Value = 6
$ws.Cells.Item(104, 18).Value = "Hortaliza"
